$d = $word.ActiveDocument

# Curly double-quote characters used in one of the defect descriptions
$lq = [char]0x201C
$rq = [char]0x201D

# Fill in the previously-blank paragraphs in the PSP0 Project Plan tables.
# Cells are addressed via Rows.Item(r).Cells.Item(c) (physical cell order)
# rather than Table.Cell(r,c), because several rows use <w:gridSpan> and
# Table.Cell() indexes by logical grid column, not physical cell position.

$t = $d.Tables.Item(2)
$t.Rows.Item(4).Cells.Item(4).Range.Text = '0,91h'
$t.Rows.Item(4).Cells.Item(6).Range.Text = '0,91h'
$t.Rows.Item(5).Cells.Item(8).Range.Text = '0'
$t.Rows.Item(11).Cells.Item(4).Range.Text = '0'
$t.Rows.Item(11).Cells.Item(6).Range.Text = '0'
$t.Rows.Item(11).Cells.Item(8).Range.Text = '0'
$t.Rows.Item(12).Cells.Item(4).Range.Text = '0'
$t.Rows.Item(12).Cells.Item(6).Range.Text = '0'
$t.Rows.Item(12).Cells.Item(8).Range.Text = '0'
$t.Rows.Item(13).Cells.Item(4).Range.Text = '2'
$t.Rows.Item(13).Cells.Item(6).Range.Text = '2'
$t.Rows.Item(19).Cells.Item(4).Range.Text = '0'
$t.Rows.Item(19).Cells.Item(6).Range.Text = '0'
$t.Rows.Item(19).Cells.Item(8).Range.Text = '0'
$t.Rows.Item(20).Cells.Item(4).Range.Text = '0'
$t.Rows.Item(20).Cells.Item(6).Range.Text = '0'
$t.Rows.Item(20).Cells.Item(8).Range.Text = '0'
$t.Rows.Item(21).Cells.Item(4).Range.Text = '2'
$t.Rows.Item(21).Cells.Item(6).Range.Text = '2'

$t = $d.Tables.Item(3)
$t.Rows.Item(4).Cells.Item(1).Range.Text = '1A'
$t.Rows.Item(4).Cells.Item(2).Range.Text = 'Code'
$t.Rows.Item(4).Cells.Item(3).Range.Text = '10/Apr/20 19:00'
$t.Rows.Item(4).Cells.Item(4).Range.Text = '0,17h'
$t.Rows.Item(4).Cells.Item(5).Range.Text = '10/Apr/20 20:05'
$t.Rows.Item(4).Cells.Item(6).Range.Text = '0,91h'
$t.Rows.Item(4).Cells.Item(7).Range.Text = 'Time to check the messages in my iPhone'

$t = $d.Tables.Item(4)
$t.Rows.Item(2).Cells.Item(2).Range.Text = '1A'
$t.Rows.Item(2).Cells.Item(4).Range.Text = '10/Apr/20'
$t.Rows.Item(2).Cells.Item(6).Range.Text = '1'
$t.Rows.Item(2).Cells.Item(8).Range.Text = '90'
$t.Rows.Item(2).Cells.Item(10).Range.Text = '1'
$t.Rows.Item(2).Cells.Item(12).Range.Text = '1'
$t.Rows.Item(2).Cells.Item(14).Range.Text = '0,08h'
$t.Rows.Item(2).Cells.Item(16).Range.Text = 'X'
$t.Rows.Item(3).Cells.Item(2).Range.Text = 'Application was not finding the right route when called by API. The configuration of the '
$t.Rows.Item(4).Cells.Item(1).Range.Text = ($lq + 'consign' + $rq + ' module was wrong.')
$t.Rows.Item(7).Cells.Item(2).Range.Text = '1A'
$t.Rows.Item(7).Cells.Item(4).Range.Text = '10/Apr/20'
$t.Rows.Item(7).Cells.Item(6).Range.Text = '2'
$t.Rows.Item(7).Cells.Item(8).Range.Text = '80'
$t.Rows.Item(7).Cells.Item(10).Range.Text = '1'
$t.Rows.Item(7).Cells.Item(12).Range.Text = '1'
$t.Rows.Item(7).Cells.Item(14).Range.Text = '0,17h'
$t.Rows.Item(7).Cells.Item(16).Range.Text = 'X'
$t.Rows.Item(8).Cells.Item(2).Range.Text = 'The result of standard deviation was wrong. I made a mistake in the code.'
